# Auto-generated Excel COM-interop script
# Applies scheduled-runner market data refresh to Chocobo_Profits sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 79.33
$ws.Range("I15").Value = 79.33
$ws.Range("K15").Value = 237.99
$ws.Range("M15").Value = -68.99000000000001
$ws.Range("H100").Value = 16667894
$ws.Range("I100").Value = 18183076
$ws.Range("J100").Value = 900
$ws.Range("K100").Value = 18183076
$ws.Range("L100").Value = 900
$ws.Range("M100").Value = -18182535
$ws.Range("N100").Value = -1982
$ws.Range("H111").Value = 2009.75
$ws.Range("I111").Value = 2009.3334
$ws.Range("J111").Value = 2011
$ws.Range("K111").Value = 6028.0002
$ws.Range("L111").Value = 6033
$ws.Range("M111").Value = -2961.0002
$ws.Range("N111").Value = -12167
$ws.Range("H137").Value = 1877.8793
$ws.Range("I137").Value = 1295.1459
$ws.Range("J137").Value = 4675
$ws.Range("K137").Value = 3885.4377
$ws.Range("L137").Value = 14025
$ws.Range("M137").Value = -1335.4377
$ws.Range("N137").Value = -19125
$ws.Range("H138").Value = 2367.7976
$ws.Range("I138").Value = 1343.6666
$ws.Range("J138").Value = 2647.106
$ws.Range("K138").Value = 4030.9998
$ws.Range("L138").Value = 7941.318000000001
$ws.Range("M138").Value = 1109.0002
$ws.Range("N138").Value = -18221.318

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H26").Value = 20250
$ws.Range("J26").Value = 35500
$ws.Range("L26").Value = 35500
$ws.Range("N26").Value = -36160
$ws.Range("H32").Value = 8690.280000000001
$ws.Range("I32").Value = 5443.3975
$ws.Range("J32").Value = 17468.889
$ws.Range("K32").Value = 5443.3975
$ws.Range("L32").Value = 17468.889
$ws.Range("M32").Value = -5156.3975
$ws.Range("N32").Value = -18042.889
$ws.Range("H45").Value = 2068.1765
$ws.Range("I45").Value = 1171.2858
$ws.Range("J45").Value = 2696
$ws.Range("K45").Value = 1171.2858
$ws.Range("L45").Value = 2696
$ws.Range("M45").Value = -794.2858000000001
$ws.Range("N45").Value = -3450
$ws.Range("H61").Value = 1410.6875
$ws.Range("I61").Value = 923.9048
$ws.Range("K61").Value = 923.9048
$ws.Range("M61").Value = -711.9048
$ws.Range("H136").Value = 1410.6875
$ws.Range("I136").Value = 923.9048
$ws.Range("K136").Value = 2771.7144
$ws.Range("M136").Value = -221.7143999999998

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2699.39
$ws.Range("I134").Value = 1620.2559
$ws.Range("K134").Value = 4860.7677
$ws.Range("M134").Value = -2325.7677

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 5052523
$ws.Range("I16").Value = 11112100
$ws.Range("J16").Value = 2875
$ws.Range("K16").Value = 11112100
$ws.Range("L16").Value = 2875
$ws.Range("M16").Value = -11111813
$ws.Range("N16").Value = -3449
$ws.Range("H31").Value = 2897.5
$ws.Range("I31").Value = 1346.3667
$ws.Range("J31").Value = 5805.875
$ws.Range("K31").Value = 1346.3667
$ws.Range("L31").Value = 5805.875
$ws.Range("M31").Value = -1051.3667
$ws.Range("N31").Value = -6395.875
$ws.Range("H34").Value = 2897.5
$ws.Range("I34").Value = 1346.3667
$ws.Range("J34").Value = 5805.875
$ws.Range("K34").Value = 1346.3667
$ws.Range("L34").Value = 5805.875
$ws.Range("M34").Value = -1144.3667
$ws.Range("N34").Value = -6209.875
$ws.Range("H36").Value = 17170.666
$ws.Range("I36").Value = 8508
$ws.Range("K36").Value = 8508
$ws.Range("M36").Value = -8120
$ws.Range("H40").Value = 17170.666
$ws.Range("I40").Value = 8508
$ws.Range("K40").Value = 8508
$ws.Range("M40").Value = -8348
$ws.Range("H58").Value = 1972.6567
$ws.Range("I58").Value = 1660.1311
$ws.Range("K58").Value = 1660.1311
$ws.Range("M58").Value = -1457.1311
$ws.Range("H113").Value = 5052523
$ws.Range("I113").Value = 11112100
$ws.Range("J113").Value = 2875
$ws.Range("K113").Value = 11112100
$ws.Range("L113").Value = 2875
$ws.Range("M113").Value = -11109930
$ws.Range("N113").Value = -7215
$ws.Range("H122").Value = 1925.7333
$ws.Range("J122").Value = 4497.8
$ws.Range("L122").Value = 13493.4
$ws.Range("N122").Value = -18393.4
$ws.Range("H134").Value = 4477.758
$ws.Range("I134").Value = 4285.533
$ws.Range("K134").Value = 12856.599
$ws.Range("M134").Value = -10321.599
$ws.Range("H136").Value = 1972.6567
$ws.Range("I136").Value = 1660.1311
$ws.Range("K136").Value = 4980.3933
$ws.Range("M136").Value = -2430.3933

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 49.727272
$ws.Range("I12").Value = 92.28570999999999
$ws.Range("K12").Value = 276.85713
$ws.Range("M12").Value = -103.85713
$ws.Range("H114").Value = 3032
$ws.Range("J114").Value = 4000
$ws.Range("L114").Value = 12000
$ws.Range("N114").Value = -18508

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H64").Value = 33392.332
$ws.Range("J64").Value = 33392.332
$ws.Range("L64").Value = 33392.332
$ws.Range("N64").Value = -33888.332
$ws.Range("H67").Value = 33392.332
$ws.Range("J67").Value = 33392.332
$ws.Range("L67").Value = 33392.332
$ws.Range("N67").Value = -35108.332
$ws.Range("H70").Value = 6563.4346
$ws.Range("I70").Value = 5890.6
$ws.Range("J70").Value = 7825
$ws.Range("K70").Value = 5890.6
$ws.Range("L70").Value = 7825
$ws.Range("M70").Value = -5620.6
$ws.Range("N70").Value = -8365
$ws.Range("H73").Value = 6563.4346
$ws.Range("I73").Value = 5890.6
$ws.Range("J73").Value = 7825
$ws.Range("K73").Value = 5890.6
$ws.Range("L73").Value = 7825
$ws.Range("M73").Value = -4954.6
$ws.Range("N73").Value = -9697
$ws.Range("H80").Value = 22729860
$ws.Range("I80").Value = 50002110
$ws.Range("J80").Value = 2983.3333
$ws.Range("K80").Value = 50002110
$ws.Range("L80").Value = 2983.3333
$ws.Range("M80").Value = -50001112
$ws.Range("N80").Value = -4979.3333
$ws.Range("H83").Value = 22729860
$ws.Range("I83").Value = 50002110
$ws.Range("J83").Value = 2983.3333
$ws.Range("K83").Value = 250010550
$ws.Range("L83").Value = 14916.6665
$ws.Range("M83").Value = -250005558
$ws.Range("N83").Value = -24900.6665
$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").ClearContents()
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 9394.25
$ws.Range("J40").Value = 7266.6665
$ws.Range("L40").Value = 7266.6665
$ws.Range("N40").Value = -7538.6665
$ws.Range("H68").Value = 872.1892
$ws.Range("I68").Value = 743.85297
$ws.Range("K68").Value = 743.85297
$ws.Range("M68").Value = 5.147029999999972
$ws.Range("H71").Value = 872.1892
$ws.Range("I71").Value = 743.85297
$ws.Range("K71").Value = 3719.26485
$ws.Range("M71").Value = 24.73514999999998
$ws.Range("H81").Value = 63819.6
$ws.Range("J81").Value = 63819.6
$ws.Range("L81").Value = 63819.6
$ws.Range("N81").Value = -65815.60000000001
$ws.Range("H84").Value = 63819.6
$ws.Range("J84").Value = 63819.6
$ws.Range("L84").Value = 191458.8
$ws.Range("N84").Value = -201442.8
$ws.Range("H122").Value = 3817.7693
$ws.Range("I122").Value = 3008.879
$ws.Range("K122").Value = 9026.636999999999
$ws.Range("M122").Value = -6576.636999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H43").Value = 22772.25
$ws.Range("J43").Value = 29029.666
$ws.Range("L43").Value = 29029.666
$ws.Range("N43").Value = -29327.666
$ws.Range("H122").Value = 9798
$ws.Range("I122").Value = 6995
$ws.Range("J122").Value = 11666.667
$ws.Range("K122").Value = 20985
$ws.Range("L122").Value = 35000.001
$ws.Range("M122").Value = -18535
$ws.Range("N122").Value = -39900.001
$ws.Range("H132").Value = 5652973
$ws.Range("I132").Value = 4252.5
$ws.Range("K132").Value = 12757.5
$ws.Range("M132").Value = -10227.5

Write-Output "Chocobo_Profits sheets updated"
